$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0006408296065709695
$ws.Range("C2").Value = 117.745847958593
$ws.Range("D2").Value = 261.3203778131603
$ws.Range("E2").Value = 2195978.878461985
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2196357.945328586

$ws.Range("B3").Value = 1.455362044514542
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 4.358119930609447
